$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. The chain-of-responsibility paragraph used to have its last sentence
#    split in two by a (hidden) "_GoBack" bookmark:
#       "...doesn't want to propagate" + bookmark + " the event further up
#       the chain, then stop."
#    Remove that now-stray bookmark; the sentence itself already reads
#    correctly once it's gone (the bookmark will be re-created further down,
#    at the end of the newly written "Command" section, matching the diff).
# ---------------------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# ---------------------------------------------------------------------------
# 2. Fill in the new "Command" design-pattern write-up. The heading
#    "Command" already exists, followed by a single placeholder bullet
#    (text ">"). Re-purpose that bullet for the first new sentence, then
#    insert two more bulleted paragraphs (inheriting the same list/style)
#    for the remaining content.
# ---------------------------------------------------------------------------
$bullet1 = $d.Paragraphs(298)
$bullet1.Range.Text = "The command design pattern is an object which represents an instruction to perform a particular action. Contains all the information necessary for the action to be taken."

$bullet1.Range.InsertParagraphAfter()
$bullet2 = $d.Paragraphs(299)
$bullet2.Range.Text = "Normally no history of an object" + [char]0x2019 + "s state is kept, rather the state is updated as needed. However, to implement various types of functionality you need to be able to record the state of an object over time. This is could be to undo or redo certain actions in a program, to record a macro " + [char]0x2013 + " a repeatable action, etc."

$bullet2.Range.InsertParagraphAfter()
$bullet3 = $d.Paragraphs(300)
$bullet3.Range.Text = "Command and query separation refer to how commands are used to make changes to an object, while queries are used to get data from an object. The term command here has a different meaning/context to the design pattern. The design pattern can be used to represent commands and queries.  "

# Re-insert the "_GoBack" bookmark right before the two trailing spaces at
# the end of the third bullet, mirroring its new position in the diff.
$bullet3 = $d.Paragraphs(300)
$r3 = $bullet3.Range
$bmPos = $r3.End - 3
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ---------------------------------------------------------------------------
# 3. The headings that follow ("Interpreter", "Mediator", "Memento",
#    "Observer", "State") each shift forward by one slot now that "Command"
#    owns real content instead of being an empty placeholder section.
# ---------------------------------------------------------------------------
$d.Paragraphs(302).Range.Text = "Interpreter"
$d.Paragraphs(305).Range.Text = "Mediator"
$d.Paragraphs(308).Range.Text = "Memento"
$d.Paragraphs(311).Range.Text = "Observer"
$d.Paragraphs(314).Range.Text = "State"
